$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 623 entirely; this shifts all subsequent rows (624-696) up by one,
# which matches the target diff (the post "「微笑みに微笑み。先にした人の方が美しい」" was removed).
$ws.Rows.Item(623).Delete()
